# Applies the "Agregar Ejemplo y cambios paciente" revision:
#  - bump IG Version / Date on the Metadata sheet
#  - Extension.value[x] (the slicing-parent row on the Elements sheet) is no
#    longer constrained to a single fixed type ("code"); now that a
#    valueCode slice exists, its Type(s) lists every open-type choice and
#    its Slicing Rules flips from "closed" to "open"
#  - widen the Type(s) column so the long type list is visible

$wb = $excel.ActiveWorkbook

# ----- Metadata sheet -----------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Plain ".Value = '1.8.11'" gets auto-parsed by Excel's input heuristics as
# the date 2011-01-08 (same trap as "1/8/11"). Enter it as a formula
# (never date-sniffed) and then paste-special the computed value back over
# itself so the cell ends up a normal shared-string text cell, not a
# formula, with no stray helper cells left behind.
$meta.Range("B3").Formula = '="1.8.11"'
$meta.Range("B3").Copy()
$meta.Range("B3").PasteSpecial(-4163)

$meta.Range("B8").Value = "2024-06-13T17:23:26-04:00"

# ----- Elements sheet ------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

$types = "base64Binary" + "`n" + "booleancanonicalcodedatedateTimedecimalidinstantintegermarkdownoidpositiveIntstringtimeunsignedInturiurluuidAddressAgeAnnotationAttachmentCodeableConceptCodingContactPointCountDistanceDurationHumanNameIdentifierMoneyPeriodQuantityRangeRatioReferenceSampledDataSignatureTimingContactDetailContributorDataRequirementExpressionParameterDefinitionRelatedArtifactTriggerDefinitionUsageContextDosageMeta"

# Row 6 = "Extension.value[x]" slicing-parent row
$elements.Range("K6").Value = $types
$elements.Range("AE6").Value = "open"

# Widen column K (Type(s)) so the long list of types fits, matching the
# bestFit recalculation Excel performs after the text grows.
$elements.Columns.Item(11).ColumnWidth = 254.14
